$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CLIENTES")

$startCol = 58  # BF
$endCol = 74    # BV
$targetRange = $ws.Range($ws.Cells.Item(1, $startCol), $ws.Cells.Item(1, $endCol))

# Reuse the existing header formatting (fontId=3 / fillId=2, the dark band
# used across row 1) from BE1, then nudge the alignment so the new style
# entry mirrors the one the source workbook picked up
# (applyFont + applyFill + applyAlignment, all default values).
$ws.Range("BE1").Copy()
$targetRange.PasteSpecial(-4122)
$targetRange.ShrinkToFit = $False
$excel.CutCopyMode = 0

# Write the new header captions in the order they were originally typed so
# the shared-string table grows in the same sequence as the source file.
$ws.Range("BF1").Value = "Categoria del cliente"
$ws.Range("BG1").Value = "Moneda"
$ws.Range("BJ1").Value = "Uso del CFDI"
$ws.Range("BK1").Value = "Clave del Producto o Servicio"
$ws.Range("BL1").Value = "Clave de Unidad"
$ws.Range("BM1").Value = "Unidad"
$ws.Range("BQ1").Value = "Base"
$ws.Range("BN1").Value = "Impuesto"
$ws.Range("BO1").Value = "Tipo de Factor"
$ws.Range("BP1").Value = "Tasa O Cuota"
$ws.Range("BT1").Value = "Cuenta"
$ws.Range("BU1").Value = "CLABE"
$ws.Range("BS1").Value = "RFC del Banco"
$ws.Range("BR1").Value = "Banco"
$ws.Range("BV1").Value = "Correo"
$ws.Range("BH1").Value = "Forma de Pago"
$ws.Range("BI1").Value = "Método de Pago CFDI"

# Match the per-column widths the author set for the new fields as closely
# as the host's column-width quantization allows.
$ws.Range("BF1").ColumnWidth = 26.833333333333332
$ws.Range("BG1:BH1").ColumnWidth = 19
$ws.Range("BI1").ColumnWidth = 28
$ws.Range("BJ1").ColumnWidth = 16.666666666666668
$ws.Range("BK1").ColumnWidth = 37.5
$ws.Range("BL1").ColumnWidth = 21
$ws.Range("BN1").ColumnWidth = 12.666666666666666
$ws.Range("BO1").ColumnWidth = 18.166666666666668
$ws.Range("BP1").ColumnWidth = 17.666666666666668
$ws.Range("BQ1").ColumnWidth = 16.5
$ws.Range("BS1").ColumnWidth = 18.5
$ws.Range("BT1").ColumnWidth = 17.833333333333332

# Restore the view state (zoom + scroll/selection) to where the author left
# the workbook after adding the new columns.
$win = $ws.Application.ActiveWindow
$win.Zoom = 106
$ws.Range("BT2").Select() | Out-Null
